# Applies the ES1_Murder.docx edit:
#  1. Extends the final body paragraph's text with the new sentence about
#     gangs replacing support systems, quoting Spergel.
#  2. Inserts a new footnote (id 20) referencing Spergel, 'Youth Gangs', p. 171.
#  3. Refreshes the cached NUMPAGES field result (4 -> 5) in both footers
#     that carry it, since the document now spans one more page.

$d = $word.ActiveDocument

# --- 1 & 2: extend the paragraph text, then append the footnote mark ---
$old = ', a stark contrast where for the American youth violence is not influenced by movies, fashion or other media, they are violent because it is the only tool that works to keep them safe in what can be seen as a broken and unsafe world. '
$new = ', a stark contrast where for the American youth violence is not influenced by movies, fashion or other media, they are violent because it is the only tool that works to keep them safe in what can be seen as a broken and unsafe world. When the support systems were removed such as schools, gangs took over, occurring in spaces where the gang performs cultural and economic functions "no longer adequately performed by the family, the school, and the labor market"'

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence to extend."
}

# Assign .Text directly (rather than a Find/Replace substitution) so the
# straight quotes in the new clause are preserved verbatim instead of being
# auto-corrected into curly quotes, and so the edit lands inside the
# existing run (keeping its run formatting) rather than spawning a new one.
$rng.Text = $new

# Collapse to the end of the (now extended) run and drop the footnote there.
$rng.Collapse(0)
$footnoteText = 'Spergel, ‘Youth Gangs’, p. 171.'
$fn = $d.Footnotes.Add($rng, "", $footnoteText)

Write-Output "Inserted footnote #$($fn.Index), total footnotes now $($d.Footnotes.Count)"

# --- 3: bump the cached NUMPAGES field result in the footers that show it ---
# (Note: `.Exists` is unreliable for the first-page/primary footer split in
# this document, so just sweep every footer slot in every section; Find
# harmlessly reports $false on slots with no "4" to replace, e.g. the blank
# even-page footer.)
foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $footer = $sec.Footers.Item($i)
        $fr = $footer.Range
        $fr.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "5", 2) | Out-Null
    }
}

Write-Output "Done."
